$wb = $excel.ActiveWorkbook

# OFF sheet (Week 15 logged / Week 16 simulated) - Home row (row 2)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 398
$wsOff.Range("C2").Value = 298
$wsOff.Range("D2").Value = 74
$wsOff.Range("E2").Value = 34

# DEF sheet - Home row (row 2)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 388
$wsDef.Range("C2").Value = 282
$wsDef.Range("D2").Value = 91
$wsDef.Range("E2").Value = 46
$wsDef.Range("F2").Value = 6
